# Atualização de bases das ligas, do dia: 28-05-2024 às 19:13
#
# The underlying source data was re-fetched and a handful of match rows
# were re-ordered. Columns A (id) and D (Date) stay anchored to their row
# position; the rest of each row's data (B..AD: external id, league,
# teams, score, odds, ...) moves with the match it belongs to.
#
#   rows 5,6,7    -> cyclic rotation (row5<-row6, row6<-row7, row7<-row5)
#   rows 20,21    -> swap
#   rows 151,152  -> swap
#
# Re-assigning the cell values (rather than touching the shared-string
# table directly) is enough: the workbook engine re-derives / re-orders
# the shared string pool from actual cell usage on save, so team-name
# strings end up correctly deduplicated/reindexed automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

# --- rows 5, 6, 7: cyclic rotation -----------------------------------
$row5 = Get-RowValues 5
$row6 = Get-RowValues 6
$row7 = Get-RowValues 7

Set-RowValues 5 $row6
Set-RowValues 6 $row7
Set-RowValues 7 $row5

# --- rows 20, 21: swap -------------------------------------------------
$row20 = Get-RowValues 20
$row21 = Get-RowValues 21

Set-RowValues 20 $row21
Set-RowValues 21 $row20

# --- rows 151, 152: swap -----------------------------------------------
$row151 = Get-RowValues 151
$row152 = Get-RowValues 152

Set-RowValues 151 $row152
Set-RowValues 152 $row151
